$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Step 1: remove the extra gap (2 rows) between the Search section and the Login section header
$ws.Range("A27:A28").EntireRow.Delete()

# Step 2: fill the existing blank spacer row (now row 44) with new Login test case TC_LF_015
$ws.Range("B61").Copy()
$ws.Range("B44").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A44").Value = "TC_LF_015"
$ws.Range("B44").Value = "If the user is signed up for a long time, make sure the session is timed out or not."

Write-Host "Step2 done"
